$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil2")

Write-Host "D4 NumberFormat:" $ws.Range("D4").NumberFormat()
Write-Host "D10 Formula:" $ws.Range("D10").Formula()
Write-Host "H10 before Formula:" $ws.Range("H10").Formula()
Write-Host "H10 before Value:" $ws.Range("H10").Value()
